$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 previously held the text "R40"; it now holds the text "1".
$ws.Range("B11").Value = "1"
